# Scheduled runner refresh: update market-price-derived columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# across the per-class profit sheets. Some rows gain/lose the NQ or HQ
# profit cell depending on whether that side is profitable this run.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 7695.857
$ws.Range("I33").Value = 11324.667
$ws.Range("J33").Value = 1164
$ws.Range("K33").Value = 11324.667
$ws.Range("L33").Value = 1164
$ws.Range("M33").Value = -11095.667
$ws.Range("N33").Value = -1622
$ws.Range("H52").Value = 999.5
$ws.Range("I52").Value = 999.5
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 2998.5
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -2838.5
$ws.Range("N52").ClearContents()
$ws.Range("H54").Value = 4680.5
$ws.Range("I54").Value = 4680.5
$ws.Range("K54").Value = 4680.5
$ws.Range("M54").Value = -4194.5
$ws.Range("H57").Value = 66598.39999999999
$ws.Range("J57").Value = 66598.39999999999
$ws.Range("L57").Value = 199795.2
$ws.Range("N57").Value = -200793.2
$ws.Range("H86").Value = 3799.3
$ws.Range("I86").Value = 3999.75
$ws.Range("K86").Value = 3999.75
$ws.Range("M86").Value = -2876.75
$ws.Range("H89").Value = 3799.3
$ws.Range("I89").Value = 3999.75
$ws.Range("K89").Value = 19998.75
$ws.Range("M89").Value = -14382.75
$ws.Range("H112").Value = 1755.4286
$ws.Range("I112").Value = 1216.3334
$ws.Range("J112").Value = 1902.4546
$ws.Range("K112").Value = 3649.0002
$ws.Range("L112").Value = 5707.3638
$ws.Range("M112").Value = -2541.0002
$ws.Range("N112").Value = -7923.3638
$ws.Range("H135").Value = 1378.909
$ws.Range("I135").Value = 1378.909
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 12410.181
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -9875.181
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 1765.3846
$ws.Range("J137").Value = 4824
$ws.Range("L137").Value = 14472
$ws.Range("N137").Value = -19572

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2172.2222
$ws.Range("I61").Value = 2172.2222
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2172.2222
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1960.2222
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 4037.2222
$ws.Range("J74").Value = 4800
$ws.Range("L74").Value = 4800
$ws.Range("N74").Value = -6548
$ws.Range("H77").Value = 4037.2222
$ws.Range("J77").Value = 4800
$ws.Range("L77").Value = 24000
$ws.Range("N77").Value = -32736
$ws.Range("H97").Value = 13758.777
$ws.Range("I97").Value = 25799.5
$ws.Range("K97").Value = 25799.5
$ws.Range("M97").Value = -25303.5
$ws.Range("H102").Value = 2259.1428
$ws.Range("J102").Value = 3961.2
$ws.Range("L102").Value = 3961.2
$ws.Range("N102").Value = -7205.2
$ws.Range("H122").Value = 2821
$ws.Range("J122").Value = 3933.3333
$ws.Range("L122").Value = 11799.9999
$ws.Range("N122").Value = -16699.9999
$ws.Range("H136").Value = 2172.2222
$ws.Range("I136").Value = 2172.2222
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6516.6666
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3966.6666
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 99999
$ws.Range("J137").Value = 99999
$ws.Range("L137").Value = 99999
$ws.Range("N137").Value = -110199

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4002845.2
$ws.Range("I20").Value = 7144760
$ws.Range("J20").Value = 4044.818
$ws.Range("K20").Value = 7144760
$ws.Range("L20").Value = 4044.818
$ws.Range("M20").Value = -7144513
$ws.Range("N20").Value = -4538.818
$ws.Range("H107").Value = 3937.25
$ws.Range("I107").Value = 2266.1667
$ws.Range("K107").Value = 2266.1667
$ws.Range("M107").Value = -346.1667000000002
$ws.Range("H141").Value = 79994
$ws.Range("J141").Value = 79994
$ws.Range("L141").Value = 79994
$ws.Range("N141").Value = -90354

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 45000
$ws.Range("I48").Value = 45000
$ws.Range("K48").Value = 45000
$ws.Range("M48").Value = -44524
$ws.Range("H105").Value = 761.2143
$ws.Range("I105").Value = 847.36365
$ws.Range("J105").Value = 445.33334
$ws.Range("K105").Value = 847.36365
$ws.Range("L105").Value = 445.33334
$ws.Range("M105").Value = 899.63635
$ws.Range("N105").Value = -3939.33334
$ws.Range("H107").Value = 1982.3889
$ws.Range("I107").Value = 389.6875
$ws.Range("K107").Value = 389.6875
$ws.Range("M107").Value = 1530.3125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 913.3333
$ws.Range("I5").Value = 370
$ws.Range("K5").Value = 1110
$ws.Range("M5").Value = -998
$ws.Range("H34").Value = 806.7143
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H38").Value = 200
$ws.Range("J38").Value = 164.92308
$ws.Range("L38").Value = 494.76924
$ws.Range("N38").Value = -1188.76924
$ws.Range("H39").Value = 5083.1665
$ws.Range("I39").Value = 597.6
$ws.Range("J39").Value = 8287.143
$ws.Range("K39").Value = 1792.8
$ws.Range("L39").Value = 24861.429
$ws.Range("M39").Value = -1498.8
$ws.Range("N39").Value = -25449.429
$ws.Range("H40").Value = 222.22223
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 222.22223
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 888.88892
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1026.88892
$ws.Range("H50").Value = 285812.16
$ws.Range("I50").Value = 166780.83
$ws.Range("K50").Value = 500342.49
$ws.Range("M50").Value = -499861.49
$ws.Range("H53").Value = 285812.16
$ws.Range("I53").Value = 166780.83
$ws.Range("K53").Value = 500342.49
$ws.Range("M53").Value = -499861.49
$ws.Range("H55").Value = 6899.6665
$ws.Range("I55").Value = 2200.75
$ws.Range("J55").Value = 9249.125
$ws.Range("K55").Value = 6602.25
$ws.Range("L55").Value = 27747.375
$ws.Range("M55").Value = -6425.25
$ws.Range("N55").Value = -28101.375
$ws.Range("H82").Value = 7487.4165
$ws.Range("I82").Value = 5753
$ws.Range("J82").Value = 8354.625
$ws.Range("K82").Value = 17259
$ws.Range("L82").Value = 25063.875
$ws.Range("M82").Value = -16853
$ws.Range("N82").Value = -25875.875
$ws.Range("H85").Value = 7487.4165
$ws.Range("I85").Value = 5753
$ws.Range("J85").Value = 8354.625
$ws.Range("K85").Value = 17259
$ws.Range("L85").Value = 25063.875
$ws.Range("M85").Value = -15855
$ws.Range("N85").Value = -27871.875
$ws.Range("H135").Value = 913.3333
$ws.Range("I135").Value = 370
$ws.Range("K135").Value = 3330
$ws.Range("M135").Value = -795

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H63").Value = 36665.668
$ws.Range("J63").Value = 36665.668
$ws.Range("L63").Value = 36665.668
$ws.Range("N63").Value = -38037.668
$ws.Range("H66").Value = 36665.668
$ws.Range("J66").Value = 36665.668
$ws.Range("L66").Value = 109997.004
$ws.Range("N66").Value = -116861.004
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H97").Value = 897.65216
$ws.Range("J97").Value = 738.5
$ws.Range("L97").Value = 738.5
$ws.Range("N97").Value = -1730.5
$ws.Range("H122").Value = 3031.7
$ws.Range("I122").Value = 2409
$ws.Range("K122").Value = 7227
$ws.Range("M122").Value = -4777

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 10420
$ws.Range("J121").Value = 10420
$ws.Range("L121").Value = 10420
$ws.Range("N121").Value = -13914
$ws.Range("H132").Value = 84609.664
$ws.Range("I132").Value = 102845.586
$ws.Range("K132").Value = 308536.758
$ws.Range("M132").Value = -306006.758

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1134.75
$ws.Range("I100").Value = 946.4167
$ws.Range("J100").Value = 1699.75
$ws.Range("K100").Value = 1892.8334
$ws.Range("L100").Value = 3399.5
$ws.Range("M100").Value = -1351.8334
$ws.Range("N100").Value = -4481.5
$ws.Range("H129").Value = 71142.336
$ws.Range("J129").Value = 71142.336
$ws.Range("L129").Value = 71142.336
$ws.Range("N129").Value = -81142.336
